$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update confidential disclaimer date (2021-05-03 -> 2021-05-04)
$ws.Range("A42").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-39
$ws.Range("D2").Value2 = 0.06065325471595705
$ws.Range("E2").Value2 = -0.03538554398672089
$ws.Range("D3").Value2 = 0.05360778749639757
$ws.Range("E3").Value2 = -0.0161597713015168
$ws.Range("D4").Value2 = 0.3042443081368644
$ws.Range("E4").Value2 = 0.003917727717923647
$ws.Range("D5").Value2 = 0.03604030736891038
$ws.Range("E5").Value2 = -0.02203461401037654
$ws.Range("D6").Value2 = 0.03218180672674929
$ws.Range("E6").Value2 = -0.002493443961996578
$ws.Range("D7").Value2 = 0.02937807218852846
$ws.Range("E7").Value2 = 0.01382368283776714
$ws.Range("D8").Value2 = 0.02813163685310837
$ws.Range("E8").Value2 = 0.01549543005871312
$ws.Range("D9").Value2 = 0.0241999166488939
$ws.Range("E9").Value2 = -0.009850830284266809
$ws.Range("D10").Value2 = 0.02493594352558151
$ws.Range("E10").Value2 = -0.01547108933540464
$ws.Range("D11").Value2 = 0.02403112853055576
$ws.Range("E11").Value2 = -0.0130820261640523
$ws.Range("D12").Value2 = 0.02244605272063885
$ws.Range("E12").Value2 = 0.01084812623274156
$ws.Range("D13").Value2 = 0.02171694338978477
$ws.Range("E13").Value2 = -0.006792086679963338
$ws.Range("D14").Value2 = 0.02155081586599799
$ws.Range("E14").Value2 = -0.003407407407407415
$ws.Range("D15").Value2 = 0.02108914949944639
$ws.Range("E15").Value2 = 0.007569564295879072
$ws.Range("D16").Value2 = 0.02159870656792243
$ws.Range("E16").Value2 = 0.01340231584134011
$ws.Range("D17").Value2 = 0.0201885914557028
$ws.Range("E17").Value2 = -0.009198734844491185
$ws.Range("D18").Value2 = 0.01462518182503164
$ws.Range("E18").Value2 = -0.006287111421585778
$ws.Range("D19").Value2 = 0.01685710138227465
$ws.Range("E19").Value2 = 0.00106063284426372
$ws.Range("D20").Value2 = 0.01541016364146424
$ws.Range("E20").Value2 = 0.01174033149171261
$ws.Range("D21").Value2 = 0.0162756018371299
$ws.Range("E21").Value2 = 0.006290377422645488
$ws.Range("D22").Value2 = 0.01457792966579953
$ws.Range("E22").Value2 = -0.01649875894291131
$ws.Range("D23").Value2 = 0.01507471588042905
$ws.Range("E23").Value2 = -0.006240822320117467
$ws.Range("D24").Value2 = 0.01441510127925643
$ws.Range("E24").Value2 = 0.01650793650793658
$ws.Range("D25").Value2 = 0.01361426231929774
$ws.Range("E25").Value2 = -0.008989642368575357
$ws.Range("D26").Value2 = 0.01399291813584698
$ws.Range("E26").Value2 = -0.01664854011545214
$ws.Range("D27").Value2 = 0.0127418001540165
$ws.Range("E27").Value2 = -0.01112531008043294
$ws.Range("D28").Value2 = 0.01347708206422974
$ws.Range("E28").Value2 = 0.005496067468965915
$ws.Range("D29").Value2 = 0.01441212141335891
$ws.Range("E29").Value2 = 0.003012804418779824
$ws.Range("D30").Value2 = 0.0131676017060158
$ws.Range("E30").Value2 = 0.01551791025475246
$ws.Range("D31").Value2 = 0.01241241854855841
$ws.Range("E31").Value2 = -0.01255230125522999
$ws.Range("D32").Value2 = 0.01342514725858723
$ws.Range("E32").Value2 = -0.006801534705266787
$ws.Range("D33").Value2 = 0.01212911844072983
$ws.Range("E33").Value2 = 0.0001316135825215614
$ws.Range("D34").Value2 = 0.006315932193577197
$ws.Range("E34").Value2 = -0.03272279980453952
$ws.Range("D35").Value2 = 0.005418141168167029
$ws.Range("E35").Value2 = -0.01164777749405832
$ws.Range("D36").Value2 = 0.005514667538490289
$ws.Range("E36").Value2 = -0.03446678760276345
$ws.Range("D37").Value2 = 0.005358011731306344
$ws.Range("E37").Value2 = -0.0253446152623843
$ws.Range("D38").Value2 = 0.004790560125392757
$ws.Range("E38").Value2 = -0.01244057404363075
$ws.Range("E39").Value2 = -0.003468989599842209

$ws.Protect("D382")
